# Insert a new weekly price record as row 102, pushing the existing
# rows 102-201 down to 103-202 (mirrors the target diff: dimension grows
# from A1:R201 to A1:R202, and every record from the old row 102 onward
# shifts down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 102 — shifts 102:201 down to 103:202.
$ws.Rows.Item(102).Insert()

# Populate the new row 102 with the new weekly record.
$ws.Range("A102").Value = 8
$ws.Range("B102").Value = "Terminal La Palmera de La Serena"
$ws.Range("C102").Value = "Coquimbo"
$ws.Range("D102").Value = 44629
$ws.Range("E102").Value = 4
$ws.Range("F102").Value = 100112021
$ws.Range("G102").Value = "Ají"
$ws.Range("H102").Value = "Inferno"
$ws.Range("I102").Value = "Primera"
$ws.Range("J102").Value = 560
$ws.Range("K102").Value = 15000
$ws.Range("L102").Value = 16000
$ws.Range("M102").Value = 15500
$ws.Range("N102").Value = "`$/caja 15 kilos"
$ws.Range("O102").Value = "Provincia de Limarí"
$ws.Range("P102").Value = 1033
$ws.Range("Q102").Value = 15
$ws.Range("R102").Value = "Hortaliza"
